# Add main picture URL to the News INSERT VALUES statements.
$d = $word.ActiveDocument

$img = "http://cep.com.vn/uploads/images/chu-de-tieng-Anh-tin-tuc-news.png"

# Map of the literal text each target paragraph currently holds, to its
# replacement (straight single quotes must be preserved, so we assign
# Range.Text directly rather than using Find/Replace, which would let
# AutoCorrect turn them into curly quotes).
$replacements = @{
    "VALUES('N000001','Hehe test news 1','Day la test new 1 nhe anh em oi',1,'2019-01-06 21:25:33',1,NULL)," = `
        "VALUES('N000001','$img','Hehe test news 1','Day la test new 1 nhe anh em oi',1,'2019-01-06 21:25:33',1,NULL),"
    "('N000002','News for iphone 5','Day la test new 2 nhe anh em oi',2,'2019-01-06 21:25:33',1,'PM00002')," = `
        "('N000002','$img','News for iphone 5','Day la test new 2 nhe anh em oi',2,'2019-01-06 21:25:33',1,'PM00002'),"
    "('N000003','News for iphone 5s','Day la test new 3 nhe anh em oi',3,'2019-01-06 21:25:33',1,'PM00003')," = `
        "('N000003','$img','News for iphone 5s','Day la test new 3 nhe anh em oi',3,'2019-01-06 21:25:33',1,'PM00003'),"
    "('N000004','News for iphone 6','Day la test new 4 nhe anh em oi',4,'2019-01-06 21:25:33',1,'PM00004')" = `
        "('N000004','$img','News for iphone 6','Day la test new 4 nhe anh em oi',4,'2019-01-06 21:25:33',1,'PM00004')"
}

foreach ($para in $d.Paragraphs) {
    $r = $para.Range
    # Drop the trailing paragraph-mark character from the comparison/text range.
    $len = $r.End - $r.Start
    $body = $d.Range($r.Start, $r.End - 1)
    $old = $body.Text
    if ($replacements.ContainsKey($old)) {
        $body.Text = $replacements[$old]
    }
}
